$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - C2, D2 change
$ws.Range("C2").Value = 343415789371543.8
$ws.Range("D2").Value = 343415789371543.8

# Row 3 (RandomForestRegressor) - B3, C3, D3 change
$ws.Range("B3").Value = 2275332136062.558
$ws.Range("C3").Value = 817951067457.6162
$ws.Range("D3").Value = 59720165437110.03

# Row 4: label change GradientBoostingRegressor -> DecisionTreeRegressor, values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 642405352252.6123
$ws.Range("C4").Value = 858794523537.6841
$ws.Range("D4").Value = 11630917956572.65

# Row 5: label change AdaBoostRegressor -> MLPRegressor, values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 5718313024229.929
$ws.Range("C5").Value = 17710697264041.78
$ws.Range("D5").Value = 41496445988537.48
